# Set G1 to a new text value "lot size options".
# This adds a new shared string entry; existing strings referenced by
# G2:G4 ("1-21", "1-5", "set to 5") remain the same displayed values,
# Excel will manage the shared string table / indices automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G1").Value = "lot size options"
